$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("D2").Value = "22.152.90"
$ws.Range("E2").Value = "  -0.19%  "
$ws.Range("D3").Value = "1.559.98"
$ws.Range("E3").Value = "  +0.49%  "
$ws.Range("D4").Value = "'0.9954"
$ws.Range("E4").Value = "  -0.94%  "
$ws.Range("D5").Value = "'0.9988"
$ws.Range("E5").Value = "  -0.74%  "
$ws.Range("D6").Value = "'288.96"
$ws.Range("E6").Value = "  +0.80%  "
$ws.Range("D7").Value = "'0.3964"
$ws.Range("E7").Value = "  +3.94%  "
$ws.Range("D8").Value = "'0.3234"
$ws.Range("E8").Value = "  -0.70%  "
$ws.Range("D9").Value = "'42.76"
$ws.Range("E9").Value = "  -2.81%  "
$ws.Range("D10").Value = "'0.07309"
$ws.Range("E10").Value = "  -0.40%  "
$ws.Range("D11").Value = "'1.104"
$ws.Range("E11").Value = "  -1.90%  "
$ws.Range("D12").Value = "'0.9954"
$ws.Range("E12").Value = "  -1.00%  "
$ws.Range("D13").Value = "'19.15"
$ws.Range("E13").Value = "  -5.38%  "
$ws.Range("D14").Value = "'5.661"
$ws.Range("E14").Value = "  -2.71%  "
$ws.Range("B15").Value = "ShibaInu"
$ws.Range("C15").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D15").Value = "'0.00001139"
$ws.Range("E15").Value = "  +5.85%  "
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").Value = "'6.696"
$ws.Range("E16").Value = "  -0.86%  "
$ws.Range("D17").Value = "1.550.95"
$ws.Range("E17").Value = "  -1.84%  "
$ws.Range("D18").Value = "'0.06604"
$ws.Range("E18").Value = "  -2.03%  "
$ws.Range("D19").Value = "'84.09"
$ws.Range("E19").Value = "  -1.68%  "
$ws.Range("D20").Value = "'0.9999"
$ws.Range("E20").Value = "  -0.66%  "
$ws.Range("D21").Value = "'6.335"
$ws.Range("E21").Value = "  +0.09%  "
$ws.Range("D22").Value = "'15.87"
$ws.Range("E22").Value = "  -1.57%  "
$ws.Range("D23").Value = "'11.32"
$ws.Range("E23").Value = "  -2.19%  "
$ws.Range("D24").Value = "22.038.04"
$ws.Range("E24").Value = "  -0.73%  "
$ws.Range("D25").Value = "'2.366"
$ws.Range("E25").Value = "  +3.17%  "
$ws.Range("D26").Value = "'2.458"
$ws.Range("E26").Value = "  -1.46%  "
$ws.Range("D27").Value = "'148.41"
$ws.Range("E27").Value = "  -0.81%  "
$ws.Range("D28").Value = "'18.74"
$ws.Range("E28").Value = "  -3.56%  "
$ws.Range("D29").Value = "'4.877"
$ws.Range("E29").Value = "  -1.33%  "
$ws.Range("D30").Value = "1.726.31"
$ws.Range("E30").Value = "  -1.24%  "
$ws.Range("D31").Value = "'119.89"
$ws.Range("E31").Value = "  -2.40%  "
$ws.Range("D32").Value = "'1.072"
$ws.Range("E32").Value = "  +2.93%  "
$ws.Range("D33").Value = "'5.745"
$ws.Range("E33").Value = "  -2.21%  "
$ws.Range("D34").Value = "'0.08378"
$ws.Range("E34").Value = "  +1.43%  "
$ws.Range("D35").Value = "'9.279"
$ws.Range("E35").Value = "  -1.18%  "
$ws.Range("D36").Value = "'1.603"
$ws.Range("E36").Value = "  -13.65%  "
$ws.Range("D37").Value = "'0.06223"
$ws.Range("E37").Value = "  +0.18%  "
$ws.Range("D38").Value = "'0.02279"
$ws.Range("E38").Value = "  -3.09%  "
$ws.Range("D39").Value = "'5.169"
$ws.Range("E39").Value = "  -1.22%  "
$ws.Range("D40").Value = "'1.217"
$ws.Range("E40").Value = "  -3.72%  "
$ws.Range("D41").Value = "'0.2077"
$ws.Range("E41").Value = "  -3.62%  "
$ws.Range("D42").Value = "'1.000"
$ws.Range("E42").Value = "  -0.58%  "
$ws.Range("D43").Value = "'10.78"
$ws.Range("E43").Value = "  -1.61%  "
$ws.Range("D44").Value = "'0.5830"
$ws.Range("E44").Value = "  -3.06%  "
$ws.Range("D45").Value = "'13.27"
$ws.Range("E45").Value = "  -2.18%  "
$ws.Range("D46").Value = "'3.733"
$ws.Range("E46").Value = "  -0.19%  "
$ws.Range("D47").Value = "'0.5639"
$ws.Range("E47").Value = "  -3.68%  "
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").Value = "'1.906"
$ws.Range("E48").Value = "  -3.32%  "
$ws.Range("B49").Value = "Quant"
$ws.Range("C49").Value = "https://coinranking.com/coin/bauj_21eYVwso+quant-qnt"
$ws.Range("D49").Value = "'118.35"
$ws.Range("E49").Value = "  -3.70%  "
$ws.Range("D50").Value = "'1.147"
$ws.Range("E50").Value = "  -1.95%  "
$ws.Range("D51").Value = "'0.06863"
$ws.Range("E51").Value = "  -2.75%  "
